$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I5").Value = "b"
$ws.Range("J5").Value = "Acknowledge (Backchannel)"
$ws.Range("I7").Value = "sd"
$ws.Range("J7").Value = "Statement-non-opinion"
$ws.Range("I13").Value = "%"
$ws.Range("J13").Value = "Uninterpretable"
$ws.Range("I17").Value = "b"
$ws.Range("J17").Value = "Acknowledge (Backchannel)"
$ws.Range("I18").Value = "b"
$ws.Range("J18").Value = "Acknowledge (Backchannel)"
$ws.Range("I23").Value = "aa"
$ws.Range("J23").Value = "Agree/Accept"
$ws.Range("I37").Value = "b"
$ws.Range("J37").Value = "Acknowledge (Backchannel)"
$ws.Range("I43").Value = "sd"
$ws.Range("J43").Value = "Statement-non-opinion"
$ws.Range("I46").Value = "sd"
$ws.Range("J46").Value = "Statement-non-opinion"
$ws.Range("I50").Value = "b"
$ws.Range("J50").Value = "Acknowledge (Backchannel)"
$ws.Range("I53").Value = "ba"
$ws.Range("J53").Value = "Appreciation"
$ws.Range("I59").Value = "aa"
$ws.Range("J59").Value = "Agree/Accept"
$ws.Range("I66").Value = "sd"
$ws.Range("J66").Value = "Statement-non-opinion"
$ws.Range("I70").Value = "sv"
$ws.Range("J70").Value = "Statement-opinion"
$ws.Range("I74").Value = "sv"
$ws.Range("J74").Value = "Statement-opinion"
$ws.Range("I84").Value = "sd"
$ws.Range("J84").Value = "Statement-non-opinion"
$ws.Range("I89").Value = "sd"
$ws.Range("J89").Value = "Statement-non-opinion"
$ws.Range("I90").Value = "aa"
$ws.Range("J90").Value = "Agree/Accept"
$ws.Range("I99").Value = "%"
$ws.Range("J99").Value = "Uninterpretable"
$ws.Range("I119").Value = "%"
$ws.Range("J119").Value = "Uninterpretable"
$ws.Range("I136").Value = "sv"
$ws.Range("J136").Value = "Statement-opinion"
$ws.Range("I139").Value = "sv"
$ws.Range("J139").Value = "Statement-opinion"
$ws.Range("I154").Value = "b"
$ws.Range("J154").Value = "Acknowledge (Backchannel)"
$ws.Range("I165").Value = "ba"
$ws.Range("J165").Value = "Appreciation"
$ws.Range("I171").Value = "ba"
$ws.Range("J171").Value = "Appreciation"
$ws.Range("I177").Value = "sd"
$ws.Range("J177").Value = "Statement-non-opinion"
$ws.Range("I190").Value = "sv"
$ws.Range("J190").Value = "Statement-opinion"
$ws.Range("I191").Value = "b"
$ws.Range("J191").Value = "Acknowledge (Backchannel)"
$ws.Range("I204").Value = "aa"
$ws.Range("J204").Value = "Agree/Accept"
$ws.Range("I205").Value = "sd"
$ws.Range("J205").Value = "Statement-non-opinion"
$ws.Range("I207").Value = "aa"
$ws.Range("J207").Value = "Agree/Accept"
$ws.Range("I208").Value = "aa"
$ws.Range("J208").Value = "Agree/Accept"
$ws.Range("I214").Value = "sd"
$ws.Range("J214").Value = "Statement-non-opinion"
$ws.Range("I217").Value = "b"
$ws.Range("J217").Value = "Acknowledge (Backchannel)"
$ws.Range("I225").Value = "sd"
$ws.Range("J225").Value = "Statement-non-opinion"
$ws.Range("I231").Value = "b"
$ws.Range("J231").Value = "Acknowledge (Backchannel)"
$ws.Range("I233").Value = "%"
$ws.Range("J233").Value = "Uninterpretable"
$ws.Range("I237").Value = "%"
$ws.Range("J237").Value = "Uninterpretable"
$ws.Range("I242").Value = "sd"
$ws.Range("J242").Value = "Statement-non-opinion"
$ws.Range("I243").Value = "b"
$ws.Range("J243").Value = "Acknowledge (Backchannel)"
$ws.Range("I245").Value = "ba"
$ws.Range("J245").Value = "Appreciation"
$ws.Range("I255").Value = "b"
$ws.Range("J255").Value = "Acknowledge (Backchannel)"
$ws.Range("I260").Value = "aa"
$ws.Range("J260").Value = "Agree/Accept"
$ws.Range("I262").Value = "sd"
$ws.Range("J262").Value = "Statement-non-opinion"
$ws.Range("I276").Value = "b"
$ws.Range("J276").Value = "Acknowledge (Backchannel)"
$ws.Range("I287").Value = "aa"
$ws.Range("J287").Value = "Agree/Accept"
$ws.Range("I289").Value = "%"
$ws.Range("J289").Value = "Uninterpretable"
$ws.Range("I295").Value = "ba"
$ws.Range("J295").Value = "Appreciation"
$ws.Range("I306").Value = "%"
$ws.Range("J306").Value = "Uninterpretable"
$ws.Range("I307").Value = "%"
$ws.Range("J307").Value = "Uninterpretable"
$ws.Range("I309").Value = "aa"
$ws.Range("J309").Value = "Agree/Accept"
$ws.Range("I310").Value = "aa"
$ws.Range("J310").Value = "Agree/Accept"
